# Adding new information into results page, deleting unrelevant stuff,
# changing translations.
#
# The "Contracts" sheet currently has:
#   A: Contracts (header) / state names
#   B: Contract Size (header) / numeric values
#   C: Minimal Number of Contributors (header) / numeric values
#
# It needs to become:
#   A: Contracts (header) / state names            (unchanged)
#   B: Minimum Contract Size (header) / NEW values
#   C: Maximum Contract Size (header) / old "Contract Size" values (shifted)
#   D: Minimal Number of Contributors (header) / old values (shifted, one
#      value corrected)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contracts")
$ws.Activate()

# Insert a brand-new column before the old "Contract Size" column so the
# old B/C columns shift right to C/D, preserving their data.
$ws.Columns.Item(2).Insert()

# Headers
$ws.Cells.Item(1, 2).Value = "Minimum Contract Size"
$ws.Cells.Item(1, 3).Value = "Maximum Contract Size"

# New "Minimum Contract Size" data (column B)
$minSizes = @(50, 30, 50, 35, 60, 30, 30, 50, 40, 40)
for ($i = 0; $i -lt $minSizes.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $minSizes[$i]
}

# One corrected value in the (now shifted) "Minimal Number of Contributors"
# column: Washington goes from 4 to 2.
$ws.Cells.Item(2, 4).Value = 2

# Resize columns to fit their new contents.
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(4).AutoFit()

# The old "Contract Size" column (now column C, "Maximum Contract Size")
# keeps a fixed width of 14 instead of autofit - re-apply its current width
# so it is stored as an explicit (non best-fit) width.
$existingWidth = $ws.Columns.Item(3).ColumnWidth
$ws.Columns.Item(3).ColumnWidth = $existingWidth

# Page setup was touched (paper size / orientation) when the sheet was
# prepared for printing.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where it was left after entering the new data.
$ws.Range("D3").Select() | Out-Null
